# The edit swaps the two observation records currently sitting in rows 16
# and 17 of the "Artfynd" sheet: everything that identifies/describes the
# observation (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Antal, Enhet, Ost/Nord coordinates and the
# Starttid/Sluttid) moves from row 17 into row 16, and vice versa. All the
# other columns (validation status, locality, county, dates, reporter, …)
# are identical between the two rows, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 gets what used to be in row 17 -------------------------------
$ws.Range("A16").Value = 111768476
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
# "Antal" is stored as text, not a number - a leading apostrophe keeps the
# digits as a text value instead of Excel auto-converting it to numeric.
$ws.Range("I16").Value = "'25"
$ws.Range("I16").Style = "Normal"
$ws.Range("J16").Value = "plantor/tuvor"
# Row 17 had a (blank) "Kön" entry that row 16 lacked; move that blank
# placeholder over to row 16 as well so the column now exists there too.
$ws.Range("L16").Value = "'"
$ws.Range("L16").Style = "Normal"
$ws.Range("Q16").Value = 525546.5036804043
$ws.Range("R16").Value = 6727881.884716956
$ws.Range("Z16").Value = "15:21"
$ws.Range("AB16").Value = "15:21"

# --- Row 17 gets what used to be in row 16 -------------------------------
$ws.Range("A17").Value = 111768503
$ws.Range("B17").Value = 88966
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 5754
$ws.Range("F17").Value = "Gultoppig fingersvamp"
$ws.Range("G17").Value = "Ramaria testaceoflava"
$ws.Range("H17").Value = "(Bres.) Corner"
$ws.Range("I17").Value = "'20"
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = "fruktkroppar"
$ws.Range("Q17").Value = 525545.3455456314
$ws.Range("R17").Value = 6727837.787189188
$ws.Range("Z17").Value = "15:22"
$ws.Range("AB17").Value = "15:22"

# Row 16's "Kön" column was blank/absent, so row 17 no longer has one.
$ws.Range("L17").ClearContents()
